# Weekly update: insert the new week's record (row 4) and push the
# previously-existing rows down by one, matching the source feed's
# newest-first-after-the-two-fixed-rows ordering.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make room for the new row; formatting (incl. the date style on column D)
# is inherited from the row above, same as Excel's native "Insert" does.
$ws.Rows.Item(4).Insert()

$ws.Range("A4").Value = 11
$ws.Range("B4").Value = "Vega Monumental Concepción"
$ws.Range("C4").Value = "Bíobío"
$ws.Range("D4").Value = 44672
$ws.Range("E4").Value = 8
$ws.Range("F4").Value = 100112052
$ws.Range("G4").Value = "Albahaca"
$ws.Range("H4").Value = "Sin especificar"
$ws.Range("I4").Value = "Primera"
$ws.Range("J4").Value = 140
$ws.Range("K4").Value = 3000
$ws.Range("L4").Value = 3500
$ws.Range("M4").Value = 3286
$ws.Range("N4").Value = "$/docena de matas"
$ws.Range("O4").Value = "Región Metropolitana"
$ws.Range("P4").Value = 548
$ws.Range("Q4").Value = 6
$ws.Range("R4").Value = "Hortaliza"
